# Auto-generated Excel COM-interop script to apply crypto price/volume updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.719.76'
$ws.Range("E2").Value = '  -2.91%  '
$ws.Range("D3").Value = '2.094.55'
$ws.Range("E3").Value = '  -1.06%  '
$ws.Range("E4").Value = '  -0.48%  '
$ws.Range("D5").Value = '''345.25'
$ws.Range("E5").Value = '  +0.35%  '
$ws.Range("D6").Value = '''1.007'
$ws.Range("E6").Value = '  -0.49%  '
$ws.Range("D7").Value = '''0.5171'
$ws.Range("E7").Value = '  -1.64%  '
$ws.Range("D8").Value = '''0.4457'
$ws.Range("E8").Value = '  -1.37%  '
$ws.Range("D9").Value = '''0.09514'
$ws.Range("E9").Value = '  +5.36%  '
$ws.Range("D10").Value = '''51.70'
$ws.Range("E10").Value = '  -3.58%  '
$ws.Range("D11").Value = '''1.179'
$ws.Range("E11").Value = '  +0.59%  '
$ws.Range("D12").Value = '''25.48'
$ws.Range("E12").Value = '  +4.28%  '
$ws.Range("D13").Value = '2.095.96'
$ws.Range("E13").Value = '  -0.73%  '
$ws.Range("D14").Value = '''6.771'
$ws.Range("E14").Value = '  -0.60%  '
$ws.Range("D15").Value = '''8.113'
$ws.Range("E15").Value = '  +0.39%  '
$ws.Range("D16").Value = '''99.58'
$ws.Range("E16").Value = '  +0.82%  '
$ws.Range("D17").Value = '''0.00001170'
$ws.Range("E17").Value = '  -0.05%  '
$ws.Range("E18").Value = '  -0.60%  '
$ws.Range("D19").Value = '''20.99'
$ws.Range("E19").Value = '  +8.38%  '
$ws.Range("D20").Value = '''0.06670'
$ws.Range("E20").Value = '  -0.62%  '
$ws.Range("E21").Value = '  -0.50%  '
$ws.Range("D22").Value = '''6.213'
$ws.Range("E22").Value = '  -1.95%  '
$ws.Range("D23").Value = '29.803.34'
$ws.Range("E23").Value = '  -2.83%  '
$ws.Range("D24").Value = '''12.73'
$ws.Range("E24").Value = '  -0.11%  '
$ws.Range("D25").Value = '''2.315'
$ws.Range("E25").Value = '  -2.64%  '
$ws.Range("D26").Value = '2.343.76'
$ws.Range("E26").Value = '  -0.71%  '
$ws.Range("D27").Value = '''22.00'
$ws.Range("E27").Value = '  -1.89%  '
$ws.Range("D28").Value = '''163.36'
$ws.Range("E28").Value = '  -1.45%  '
$ws.Range("D29").Value = '''2.538'
$ws.Range("E29").Value = '  -0.16%  '
$ws.Range("D30").Value = '''133.20'
$ws.Range("E30").Value = '  -1.43%  '
$ws.Range("D31").Value = '''1.153'
$ws.Range("E31").Value = '  -3.65%  '
$ws.Range("D32").Value = '''0.1056'
$ws.Range("E32").Value = '  -1.58%  '
$ws.Range("D33").Value = '''1.626'
$ws.Range("E33").Value = '  -0.86%  '
$ws.Range("D34").Value = '''6.216'
$ws.Range("E34").Value = '  -2.39%  '
$ws.Range("D35").Value = '''3.940'
$ws.Range("E35").Value = '  -0.46%  '
$ws.Range("D36").Value = '''6.166'
$ws.Range("E36").Value = '  +4.78%  '
$ws.Range("D37").Value = '''10.16'
$ws.Range("E37").Value = '  -0.80%  '
$ws.Range("D38").Value = '''0.02578'
$ws.Range("E38").Value = '  -2.43%  '
$ws.Range("D39").Value = '''0.06744'
$ws.Range("E39").Value = '  -1.58%  '
$ws.Range("D40").Value = '''0.2284'
$ws.Range("E40").Value = '  -2.05%  '
$ws.Range("D41").Value = '''0.6902'
$ws.Range("E41").Value = '  +0.28%  '
$ws.Range("E42").Value = '  -1.79%  '
$ws.Range("D43").Value = '''1.283'
$ws.Range("E43").Value = '  +1.37%  '
$ws.Range("D44").Value = '''0.6683'
$ws.Range("E44").Value = '  +3.95%  '
$ws.Range("D45").Value = '''14.15'
$ws.Range("E45").Value = '  -4.90%  '
$ws.Range("D46").Value = '''2.304'
$ws.Range("E46").Value = '  -0.58%  '
$ws.Range("D47").Value = '''3.628'
$ws.Range("E47").Value = '  -2.77%  '
$ws.Range("B48").Value = 'EOS'
$ws.Range("C48").Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range("D48").Value = '''1.221'
$ws.Range("E48").Value = '  -2.73%  '
$ws.Range("B49").Value = 'BabyDogeCoin'
$ws.Range("C49").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D49").Value = '''0.00000000341'
$ws.Range("E49").Value = '  -4.90%  '
$ws.Range("D50").Value = '''82.08'
$ws.Range("E50").Value = '  -1.03%  '
$ws.Range("D51").Value = '''0.07125'
$ws.Range("E51").Value = '  -2.21%  '
